# Update metrics values in columns B:Q for rows 2-26 (all model rows share
# identical new metric values according to the diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    [double]"0.9999674344715328",
    [double]"0.9989400190852257",
    [double]"0.9999458681674095",
    [double]"0.9999755797402436",
    [double]"0.9999604494174408",
    [double]"3.039849136261971e-05",
    [double]"0.0009894456561571146",
    [double]"3.827352771391014e-05",
    [double]"1.663551008379092e-05",
    [double]"2.745451889885053e-05",
    [double]"0.000349036811621329",
    [double]"0.005513482689065026",
    [double]"1.000060120975632",
    [double]"0.005748203061621094",
    [double]"94.80223515402628",
    [double]"139.9006406741497"
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 2]
    }
}
